$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.373.63'
$ws.Range("E2").Value = '  -3.83%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.570.57'
$ws.Range("E3").Value = '  -4.31%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.28'
$ws.Range("E5").Value = '  -4.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '185.37'
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.562.68'
$ws.Range("E7").Value = '  -4.49%  '

$ws.Range("E8").Value = '  -4.19%  '

$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.669'
$ws.Range("E10").Value = '  -7.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.146'
$ws.Range("E11").Value = '  -10.32%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.69'
$ws.Range("E12").Value = '  -6.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("E13").Value = '  -13.15%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.76'
$ws.Range("E14").Value = '  -8.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.137.13'
$ws.Range("E15").Value = '  -4.53%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.566.22'
$ws.Range("E16").Value = '  -4.65%  '

$ws.Range("E17").Value = '  -0.92%  '

$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.21'
$ws.Range("E18").Value = '  -6.73%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.27'
$ws.Range("E19").Value = '  -6.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '66.158.88'
$ws.Range("E20").Value = '  -4.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.06'
$ws.Range("E21").Value = '  -7.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '395.59'
$ws.Range("E22").Value = '  -4.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.36'
$ws.Range("E23").Value = '  -6.60%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.38'
$ws.Range("E24").Value = '  -4.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.04'
$ws.Range("E25").Value = '  +0.53%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.87'
$ws.Range("E26").Value = '  -6.23%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.45'
$ws.Range("E27").Value = '  -2.96%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.04'
$ws.Range("E28").Value = '  -0.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.54'
$ws.Range("E29").Value = '  -6.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.95'
$ws.Range("E30").Value = '  -7.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '30.91'
$ws.Range("E31").Value = '  -7.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.04'
$ws.Range("E32").Value = '  -4.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.13'
$ws.Range("E33").Value = '  -5.08%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '618.30'
$ws.Range("E34").Value = '  -0.38%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '63.52'
$ws.Range("E35").Value = '  -3.60%  '

$ws.Range("E36").Value = '  -8.35%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '41.30'
$ws.Range("E37").Value = '  -7.26%  '

$ws.Range("E38").Value = '  +0.06%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.388'
$ws.Range("E39").Value = '  -4.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0₃0751'
$ws.Range("E40").Value = '  -14.78%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.131'
$ws.Range("E41").Value = '  -8.46%  '

$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.43%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.002.18'
$ws.Range("E43").Value = '  +6.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.81'
$ws.Range("E44").Value = '  -8.82%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.48'
$ws.Range("E45").Value = '  -6.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0406'
$ws.Range("E46").Value = '  -8.83%  '

$ws.Range("E47").Value = '  -6.51%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.05'
$ws.Range("E48").Value = '  -1.99%  '

$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.60'
$ws.Range("E49").Value = '  -6.97%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '138.44'
$ws.Range("E50").Value = '  -2.06%  '

$ws.Range("E51").Value = '  -0.63%  '
